# "Commit avant changement sur Logit VNF"
# Reset the second logistic-regression model's results (columns E:G,
# the "VNF" model's OR / IC / p columns) back to the placeholder "-"
# for every variable row, ahead of the upcoming Logit VNF change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3:G44").Value2 = "-"
